$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = [double]"22.6200000000001"
$ws.Range("H2").Value = [double]"0.2148609219202342"
$ws.Range("I2").Value = [double]"0.2148609219202342"
$ws.Range("L2").Value = [double]"6.409153644773141"
$ws.Range("M2").Value = "[-3.1775345670923985, 15.99584185663868]"
$ws.Range("N2").Value = [double]"0.1848772208142804"
$ws.Range("O2").Value = [double]"0.1848772208142804"
$ws.Range("P2").Value = [double]"-0.9937370155499243"
$ws.Range("Q2").Value = "[-4.10702703262089, 2.1195530015210418]"
$ws.Range("R2").Value = [double]"0.5235608107259564"
$ws.Range("S2").Value = [double]"0.5235608107259564"
$ws.Range("T2").Value = [double]"15.41130532158703"
$ws.Range("U2").Value = "[10.314964314850542, 20.507646328323517]"
$ws.Range("V2").Value = [double]"2.290247180791738e-07"
$ws.Range("W2").Value = [double]"2.290247180791738e-07"
$ws.Range("X2").Value = [double]"3.577537537537552"
$ws.Range("Y2").Value = [double]"-7.630570570570603"
$ws.Range("Z2").Value = [double]"14.78564564564571"

# Row 3
$ws.Range("F3").Value = [double]"22.6200000000001"
$ws.Range("H3").Value = [double]"0.258415300343331"
$ws.Range("I3").Value = [double]"0.258415300343331"
$ws.Range("L3").Value = [double]"5.320904537985973"
$ws.Range("M3").Value = "[-3.5200198676590535, 14.161828943630999]"
$ws.Range("N3").Value = [double]"0.2317670809567127"
$ws.Range("O3").Value = [double]"0.2317670809567127"
$ws.Range("P3").Value = [double]"-0.754736973835386"
$ws.Range("Q3").Value = "[-3.849158566560467, 2.339684618889695]"
$ws.Range("R3").Value = [double]"0.6256408294446043"
$ws.Range("S3").Value = [double]"0.6256408294446043"
$ws.Range("T3").Value = [double]"14.13011239186586"
$ws.Range("U3").Value = "[9.55370824285535, 18.706516540876375]"
$ws.Range("V3").Value = [double]"1.477272517647066e-07"
$ws.Range("W3").Value = [double]"1.477272517647066e-07"
$ws.Range("X3").Value = [double]"2.71711711711713"
$ws.Range("Y3").Value = [double]"-8.423063063063097"
$ws.Range("Z3").Value = [double]"13.85729729729736"

# Row 4
$ws.Range("F4").Value = [double]"22.6200000000001"
$ws.Range("H4").Value = [double]"0.0506349787009861"
$ws.Range("I4").Value = [double]"0.0506349787009861"
$ws.Range("L4").Value = [double]"7.927040360378845"
$ws.Range("M4").Value = "[-0.02174796539780033, 15.87582868615549]"
$ws.Range("N4").Value = [double]"0.05060371436334443"
$ws.Range("O4").Value = [double]"0.05060371436334443"
$ws.Range("P4").Value = [double]"-1.396263401595464"
$ws.Range("Q4").Value = "[-2.717053105807389, -0.0754736973835386]"
$ws.Range("R4").Value = [double]"0.03874446860657943"
$ws.Range("S4").Value = [double]"0.03874446860657943"
$ws.Range("T4").Value = [double]"14.16167439465498"
$ws.Range("U4").Value = "[9.786855587653367, 18.536493201656587]"
$ws.Range("V4").Value = [double]"5.267329705738177e-08"
$ws.Range("W4").Value = [double]"5.267329705738177e-08"
$ws.Range("X4").Value = [double]"5.026666666666689"
$ws.Range("Y4").Value = [double]"0.2717117117117125"
$ws.Range("Z4").Value = [double]"9.781621621621664"

# Row 5
$ws.Range("F5").Value = [double]"22.6200000000001"
$ws.Range("H5").Value = [double]"0.7846418797697239"
$ws.Range("I5").Value = [double]"0.7846418797697239"
$ws.Range("L5").Value = [double]"2.357001162125802"
$ws.Range("M5").Value = "[-6.16113317516773, 10.875135499419333]"
$ws.Range("N5").Value = [double]"0.5800771814922996"
$ws.Range("O5").Value = [double]"0.5800771814922996"
$ws.Range("P5").Value = [double]"1.666710817219811"
$ws.Range("Q5").Value = "[-1.4654476241970391, 4.798869258636661]"
$ws.Range("R5").Value = [double]"0.2895415006365103"
$ws.Range("S5").Value = [double]"0.2895415006365103"
$ws.Range("T5").Value = [double]"13.52960814346943"
$ws.Range("U5").Value = "[8.832870957392004, 18.22634532954686]"
$ws.Range("V5").Value = [double]"6.140158450573807e-07"
$ws.Range("W5").Value = [double]"6.140158450573807e-07"
$ws.Range("X5").Value = [double]"16.61969969969977"
$ws.Range("Y5").Value = [double]"5.343663663663685"
$ws.Range("Z5").Value = [double]"27.89573573573585"

# Row 6
$ws.Range("F6").Value = [double]"22.6200000000001"
$ws.Range("H6").Value = [double]"0.03503585089831429"
$ws.Range("I6").Value = [double]"0.03503585089831429"
$ws.Range("L6").Value = [double]"7.707824831530122"
$ws.Range("M6").Value = "[-0.30139941519245994, 15.717049078252703]"
$ws.Range("N6").Value = [double]"0.05887402864394908"
$ws.Range("O6").Value = [double]"0.05887402864394908"
$ws.Range("P6").Value = [double]"0.5723422051585008"
$ws.Range("Q6").Value = "[-0.9434212172942322, 2.0881056276112337]"
$ws.Range("R6").Value = [double]"0.4509151626829582"
$ws.Range("S6").Value = [double]"0.4509151626829582"
$ws.Range("T6").Value = [double]"13.45409552387514"
$ws.Range("U6").Value = "[9.349958133281234, 17.558232914469045]"
$ws.Range("V6").Value = [double]"3.967235229218602e-08"
$ws.Range("W6").Value = [double]"3.967235229218602e-08"
$ws.Range("X6").Value = [double]"20.55951951951961"
$ws.Range("Y6").Value = [double]"15.10264264264271"
$ws.Range("Z6").Value = [double]"26.01639639639651"

# Row 7
$ws.Range("F7").Value = [double]"22.6200000000001"
$ws.Range("H7").Value = [double]"0.1769795655090065"
$ws.Range("I7").Value = [double]"0.1769795655090065"
$ws.Range("L7").Value = [double]"6.314283738146024"
$ws.Range("M7").Value = "[-2.7577886406382763, 15.386356116930324]"
$ws.Range("N7").Value = [double]"0.1678224528409349"
$ws.Range("O7").Value = [double]"0.1678224528409349"
$ws.Range("P7").Value = [double]"0.3836579616996545"
$ws.Range("Q7").Value = "[-2.654158357987773, 3.421474281387082]"
$ws.Range("R7").Value = [double]"0.8003691774618764"
$ws.Range("S7").Value = [double]"0.8003691774618764"
$ws.Range("T7").Value = [double]"14.9333619141991"
$ws.Range("U7").Value = "[10.122877003518703, 19.743846824879505]"
$ws.Range("V7").Value = [double]"1.31625123467316e-07"
$ws.Range("W7").Value = [double]"1.31625123467316e-07"
$ws.Range("X7").Value = [double]"21.23879879879889"
$ws.Range("Y7").Value = [double]"10.30240240240245"
$ws.Range("Z7").Value = [double]"32.17519519519533"

# Row 8
$ws.Range("F8").Value = [double]"22.6200000000001"
$ws.Range("H8").Value = [double]"0.203768491381227"
$ws.Range("I8").Value = [double]"0.203768491381227"
$ws.Range("L8").Value = [double]"6.360179662590364"
$ws.Range("M8").Value = "[-3.248423639578257, 15.968782964758987]"
$ws.Range("N8").Value = [double]"0.1891808675573059"
$ws.Range("O8").Value = [double]"0.1891808675573059"
$ws.Range("P8").Value = [double]"0.5220264069028095"
$ws.Range("Q8").Value = "[-2.5535267614763884, 3.5975795752820074]"
$ws.Range("R8").Value = [double]"0.7340460810806304"
$ws.Range("S8").Value = [double]"0.7340460810806304"
$ws.Range("T8").Value = [double]"12.93739546963193"
$ws.Range("U8").Value = "[7.906327965858038, 17.968462973405813]"
$ws.Range("V8").Value = [double]"5.052816923756254e-06"
$ws.Range("W8").Value = [double]"5.052816923756254e-06"
$ws.Range("X8").Value = [double]"20.74066066066074"
$ws.Range("Y8").Value = [double]"9.668408408408441"
$ws.Range("Z8").Value = [double]"31.81291291291305"

# Row 9
$ws.Range("F9").Value = [double]"22.6200000000001"
$ws.Range("H9").Value = [double]"0.7051422665258222"
$ws.Range("I9").Value = [double]"0.7051422665258222"
$ws.Range("L9").Value = [double]"3.092530723721935"
$ws.Range("M9").Value = "[-6.984971012847667, 13.170032460291537]"
$ws.Range("N9").Value = [double]"0.5396400587557189"
$ws.Range("O9").Value = [double]"0.5396400587557189"
$ws.Range("P9").Value = [double]"2.295658295415965"
$ws.Range("Q9").Value = "[-0.8427896207828471, 5.434106211614778]"
$ws.Range("R9").Value = [double]"0.1476467897293976"
$ws.Range("S9").Value = [double]"0.1476467897293976"
$ws.Range("T9").Value = [double]"13.81989863834949"
$ws.Range("U9").Value = "[8.54822899455398, 19.09156828214499]"
$ws.Range("V9").Value = [double]"3.601314294998659e-06"
$ws.Range("W9").Value = [double]"3.601314294998659e-06"
$ws.Range("X9").Value = [double]"14.3554354354355"
$ws.Range("Y9").Value = [double]"3.056756756756768"
$ws.Range("Z9").Value = [double]"25.65411411411422"

# Row 10
$ws.Range("F10").Value = [double]"22.6200000000001"
$ws.Range("H10").Value = [double]"0.6465905276051531"
$ws.Range("I10").Value = [double]"0.6465905276051531"
$ws.Range("L10").Value = [double]"3.70192673852893"
$ws.Range("M10").Value = "[-6.120575117051131, 13.52442859410899]"
$ws.Range("N10").Value = [double]"0.4517628006588927"
$ws.Range("O10").Value = [double]"0.4517628006588927"
$ws.Range("P10").Value = [double]"1.817658211986886"
$ws.Range("Q10").Value = "[-1.3019212798660416, 4.9372377038398145]"
$ws.Range("R10").Value = [double]"0.2467533478745199"
$ws.Range("S10").Value = [double]"0.2467533478745199"
$ws.Range("T10").Value = [double]"15.82580028228045"
$ws.Range("U10").Value = "[10.31388764380658, 21.337712920754328]"
$ws.Range("V10").Value = [double]"6.55127023296842e-07"
$ws.Range("W10").Value = [double]"6.55127023296842e-07"
$ws.Range("X10").Value = [double]"16.07627627627635"
$ws.Range("Y10").Value = [double]"4.845525525525549"
$ws.Range("Z10").Value = [double]"27.30702702702715"

# Row 11
$ws.Range("F11").Value = [double]"23.01000000000016"
$ws.Range("H11").Value = [double]"0.3628206422094586"
$ws.Range("I11").Value = [double]"0.3628206422094586"
$ws.Range("L11").Value = [double]"5.607164019937832"
$ws.Range("M11").Value = "[-4.757056235698007, 15.971384275573671]"
$ws.Range("N11").Value = [double]"0.2816667620011801"
$ws.Range("O11").Value = [double]"0.2816667620011801"
$ws.Range("P11").Value = [double]"1.981184556317888"
$ws.Range("Q11").Value = "[-1.1383949355350396, 5.100764048170815]"
$ws.Range("R11").Value = [double]"0.2074122186176206"
$ws.Range("S11").Value = [double]"0.2074122186176206"
$ws.Range("T11").Value = [double]"16.27258251211958"
$ws.Range("U11").Value = "[10.794757569904789, 21.750407454334365]"
$ws.Range("V11").Value = [double]"3.307439926381051e-07"
$ws.Range("W11").Value = [double]"3.307439926381051e-07"
$ws.Range("X11").Value = [double]"15.7545945945947"
$ws.Range("Y11").Value = [double]"4.330210210210243"
$ws.Range("Z11").Value = [double]"27.17897897897916"

# Row 12
$ws.Range("F12").Value = [double]"23.01000000000016"
$ws.Range("H12").Value = [double]"0.1077629777844736"
$ws.Range("I12").Value = [double]"0.1077629777844736"
$ws.Range("L12").Value = [double]"6.477213212605566"
$ws.Range("M12").Value = "[-1.4236000685286267, 14.378026493739759]"
$ws.Range("N12").Value = [double]"0.1056617142010823"
$ws.Range("O12").Value = [double]"0.1056617142010823"
$ws.Range("P12").Value = [double]"2.622710984077965"
$ws.Range("Q12").Value = "[-0.08176317216550011, 5.327185140321431]"
$ws.Range("R12").Value = [double]"0.05702947020819904"
$ws.Range("S12").Value = [double]"0.05702947020819904"
$ws.Range("T12").Value = [double]"14.3021340317713"
$ws.Range("U12").Value = "[9.985401100303458, 18.61886696323915]"
$ws.Range("V12").Value = [double]"3.116317448537131e-08"
$ws.Range("W12").Value = [double]"3.116317448537131e-08"
$ws.Range("X12").Value = [double]"13.40522522522532"
$ws.Range("Y12").Value = [double]"3.501021021021044"
$ws.Range("Z12").Value = [double]"23.30942942942959"

# Row 13
$ws.Range("F13").Value = [double]"23.01000000000016"
$ws.Range("H13").Value = [double]"0.239187007930773"
$ws.Range("I13").Value = [double]"0.239187007930773"
$ws.Range("L13").Value = [double]"6.007859096924089"
$ws.Range("M13").Value = "[-3.5010189949573958, 15.516737188805575]"
$ws.Range("N13").Value = [double]"0.209716543532215"
$ws.Range("O13").Value = [double]"0.209716543532215"
$ws.Range("P13").Value = [double]"2.647868883205811"
$ws.Range("Q13").Value = "[-0.48428955821104047, 5.7800273246226626]"
$ws.Range("R13").Value = [double]"0.09552586427299814"
$ws.Range("S13").Value = [double]"0.09552586427299814"
$ws.Range("T13").Value = [double]"12.26657813660204"
$ws.Range("U13").Value = "[7.224693496697579, 17.30846277650649]"
$ws.Range("V13").Value = [double]"1.281635317740282e-05"
$ws.Range("W13").Value = [double]"1.281635317740282e-05"
$ws.Range("X13").Value = [double]"13.31309309309319"
$ws.Range("Y13").Value = [double]"1.842642642642653"
$ws.Range("Z13").Value = [double]"24.78354354354372"

# Row 14
$ws.Range("F14").Value = [double]"23.01000000000016"
$ws.Range("H14").Value = [double]"0.1252181253601431"
$ws.Range("I14").Value = [double]"0.1252181253601431"
$ws.Range("L14").Value = [double]"5.498566145759433"
$ws.Range("M14").Value = "[-1.6769768795660696, 12.674109171084936]"
$ws.Range("N14").Value = [double]"0.129739444280808"
$ws.Range("O14").Value = [double]"0.129739444280808"
$ws.Range("P14").Value = [double]"-2.817684702318773"
$ws.Range("Q14").Value = "[-4.717106086471162, -0.9182633181663848]"
$ws.Range("R14").Value = [double]"0.004537861679516197"
$ws.Range("S14").Value = [double]"0.004537861679516197"
$ws.Range("T14").Value = [double]"11.75142351455216"
$ws.Range("U14").Value = "[7.948042809730682, 15.554804219373635]"
$ws.Range("V14").Value = [double]"1.455736715616496e-07"
$ws.Range("W14").Value = [double]"1.455736715616496e-07"
$ws.Range("X14").Value = [double]"10.31879879879887"
$ws.Range("Y14").Value = [double]"3.36282282282284"
$ws.Range("Z14").Value = [double]"17.2747747747749"

# Row 15
$ws.Range("F15").Value = [double]"23.01000000000016"
$ws.Range("H15").Value = [double]"0.6289114102124211"
$ws.Range("I15").Value = [double]"0.6289114102124211"
$ws.Range("L15").Value = [double]"3.136733353761305"
$ws.Range("M15").Value = "[-5.717127207559279, 11.990593915081888]"
$ws.Range("N15").Value = [double]"0.4791877437768848"
$ws.Range("O15").Value = [double]"0.4791877437768848"
$ws.Range("P15").Value = [double]"2.547237286694427"
$ws.Range("Q15").Value = "[-0.5912106295043857, 5.685685202893239]"
$ws.Range("R15").Value = [double]"0.1090907107368142"
$ws.Range("S15").Value = [double]"0.1090907107368142"
$ws.Range("T15").Value = [double]"11.75753450013903"
$ws.Range("U15").Value = "[7.0967918495699855, 16.418277150708068]"
$ws.Range("V15").Value = [double]"7.022857882788003e-06"
$ws.Range("W15").Value = [double]"7.022857882788003e-06"
$ws.Range("X15").Value = [double]"13.68162162162172"
$ws.Range("Y15").Value = [double]"2.188138138138152"
$ws.Range("Z15").Value = [double]"25.17510510510528"

# Row 16
$ws.Range("F16").Value = [double]"23.01000000000016"
$ws.Range("H16").Value = [double]"0.3405914513260067"
$ws.Range("I16").Value = [double]"0.3405914513260067"
$ws.Range("L16").Value = [double]"4.644687258671542"
$ws.Range("M16").Value = "[-3.967654989077495, 13.25702950642058]"
$ws.Range("N16").Value = [double]"0.2831670945195301"
$ws.Range("O16").Value = [double]"0.2831670945195301"
$ws.Range("P16").Value = [double]"-2.490632013656773"
$ws.Range("Q16").Value = "[-5.597632555945777, 0.6163685286322318]"
$ws.Range("R16").Value = [double]"0.1134004681979883"
$ws.Range("S16").Value = [double]"0.1134004681979883"
$ws.Range("T16").Value = [double]"13.74799747469834"
$ws.Range("U16").Value = "[9.28179482418141, 18.21420012521527]"
$ws.Range("V16").Value = [double]"1.575892969452042e-07"
$ws.Range("W16").Value = [double]"1.575892969452042e-07"
$ws.Range("X16").Value = [double]"9.121081081081142"
$ws.Range("Y16").Value = [double]"-2.257237237237254"
$ws.Range("Z16").Value = [double]"20.49939939939954"

# Row 17
$ws.Range("B17").Value = [double]"0"
$ws.Range("F17").Value = [double]"23.01000000000016"
$ws.Range("H17").Value = [double]"0.8592735945851597"
$ws.Range("I17").Value = [double]"0.8592735945851597"
$ws.Range("L17").Value = [double]"1.589024367894346"
$ws.Range("M17").Value = "[-6.066789151376604, 9.244837887165296]"
$ws.Range("N17").Value = [double]"0.6779038869595992"
$ws.Range("O17").Value = [double]"0.6779038869595992"
$ws.Range("P17").Value = [double]"-2.226474072814388"
$ws.Range("Q17").Value = "[-5.3649219890132, 0.9119738433844242]"
$ws.Range("R17").Value = [double]"0.159954486883519"
$ws.Range("S17").Value = [double]"0.159954486883519"
$ws.Range("T17").Value = [double]"13.18184472539051"
$ws.Range("U17").Value = "[9.116113670478711, 17.247575780302302]"
$ws.Range("V17").Value = [double]"5.085400878179769e-08"
$ws.Range("W17").Value = [double]"5.085400878179769e-08"
$ws.Range("X17").Value = [double]"8.15369369369375"
$ws.Range("Y17").Value = [double]"-3.339789789789812"
$ws.Range("Z17").Value = [double]"19.64717717717731"
